$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper now also reports "height" and "weight". Insert two blank
# columns right after the current data (F:G); this pushes nothing except
# shifts what will become the "fantasy points" column out of the way and
# keeps the header formatting (bold/border/centered) that already lives on
# row 1 of the neighbouring columns.
$ws.Columns("F:G").Insert()

# The existing "fantasy points" column (header + all 16 data rows) moves
# from E to G.
$ws.Range("E1:E17").Cut()
$ws.Range("G1").Select() | Out-Null
$ws.Paste()

# Column E becomes "height": same constant value (6.5) for every player row.
$ws.Range("E1").Value = "height"
$ws.Range("E2:E17").Value = 6.5

# Column F becomes "weight": same constant value (265) for every player row.
$ws.Range("F1").Value = "weight"
$ws.Range("F2:F17").Value = 265

$ws.Range("A1").Select() | Out-Null
